# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the consolidated "全部类型" sheet to match the latest scrape.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 13230
$wsExhibit.Range("F6").Value = 111
$wsExhibit.Range("F9").Value = 38
$wsExhibit.Range("F11").Value = 13168
$wsExhibit.Range("F15").Value = 7905
$wsExhibit.Range("F17").Value = 136
$wsExhibit.Range("F27").Value = 83

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 13230
$wsAll.Range("F7").Value = 111
$wsAll.Range("F10").Value = 38
$wsAll.Range("F12").Value = 13168
$wsAll.Range("F16").Value = 7905
$wsAll.Range("F18").Value = 136
$wsAll.Range("F30").Value = 83
